$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Unit column values from "EJ/y" to "EJ/yr"
$ws.Range("D2").Value = "EJ/yr"
$ws.Range("D3").Value = "EJ/yr"

# Update the selection to D3 (matches the sheetView selection in the diff)
$ws.Range("D3").Select()
